{"js": "// Remove everything in the body after the first (Title) paragraph:\n// the UNCLASSIFIED / USA / KMF_Memo header block, the blank spacer\n// paragraph, the BLUF summary paragraph, the [Analyst Comment]\n// paragraph, and the trailing blank paragraph - leaving only the\n// Title paragraph before the section properties.\n\nconst body = context.document.body;\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\n// Remember the first (Title) paragraph's style so it survives the\n// paragraph-mark merge performed below.\nconst titleStyle = paragraphs.items[0].style;\n\nconst count = paragraphs.items.length;\n\n// Delete every paragraph except the first one. Word (and Office.js)\n// never lets you remove the body's final paragraph outright, so the\n// very last paragraph will remain behind as an empty stub after this\n// loop - it gets folded into the title paragraph next.\nfor (let i = count - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\nparagraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nif (paragraphs.items.length > 1) {\n  // Merge the leftover trailing paragraph into the title paragraph by\n  // deleting the paragraph mark between them.\n  const mergeStart = paragraphs.items[0].getRange(\"End\");\n  const mergeEnd = paragraphs.items[1].getRange(\"Start\");\n  mergeStart.expandTo(mergeEnd).delete();\n  await context.sync();\n}\n\n// The merge adopts the trailing (non-title) paragraph's formatting, so\n// restore the Title style on the single remaining paragraph.\nparagraphs = body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\nparagraphs.items[0].style = titleStyle;\nawait context.sync();\n", "ps1": "# Remove everything in the body after the first (Title) paragraph:\n# the UNCLASSIFIED / USA / KMF_Memo header block, the blank spacer\n# paragraph, the BLUF summary paragraph, the [Analyst Comment]\n# paragraph, and the trailing blank paragraph - leaving only the\n# Title paragraph before the section properties.\n\n$d = $word.ActiveDocument\n\n# Remember the first (Title) paragraph's style so it survives the\n# paragraph-mark merge performed below.\n$titleStyle = $d.Paragraphs.Item(1).Range.Style\n\n# Delete every paragraph except the first one. Word never lets you\n# remove the body's final paragraph outright, so the very last\n# paragraph will remain behind as an empty stub after this loop - it\n# gets folded into the title paragraph next.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 2; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\nif ($d.Paragraphs.Count -gt 1) {\n    # Merge the leftover trailing paragraph into the title paragraph by\n    # deleting the paragraph mark between them.\n    $firstPara = $d.Paragraphs.Item(1)\n    $mergeRange = $d.Range($firstPara.Range.End - 1, $firstPara.Range.End)\n    $mergeRange.Delete()\n}\n\n# The merge adopts the trailing (non-title) paragraph's formatting, so\n# restore the Title style on the single remaining paragraph.\n$d.Paragraphs.Item(1).Range.set_Style($titleStyle)\n"}
